# Auto-generated script applying Shinryu_Profits value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2847.9678
$ws.Range("I62").Value = 2653.9443
$ws.Range("J62").Value = 3116.6155
$ws.Range("K62").Value = 2653.9443
$ws.Range("L62").Value = 3116.6155
$ws.Range("M62").Value = -2029.9443
$ws.Range("N62").Value = -4364.6155
$ws.Range("H64").Value = 3554.2126
$ws.Range("I64").Value = 3492
$ws.Range("J64").Value = 3624.9092
$ws.Range("K64").Value = 3492
$ws.Range("L64").Value = 3624.9092
$ws.Range("M64").Value = -3244
$ws.Range("N64").Value = -4120.9092
$ws.Range("H65").Value = 2847.9678
$ws.Range("I65").Value = 2653.9443
$ws.Range("J65").Value = 3116.6155
$ws.Range("K65").Value = 13269.7215
$ws.Range("L65").Value = 15583.0775
$ws.Range("M65").Value = -10149.7215
$ws.Range("N65").Value = -21823.0775
$ws.Range("H67").Value = 3554.2126
$ws.Range("I67").Value = 3492
$ws.Range("J67").Value = 3624.9092
$ws.Range("K67").Value = 3492
$ws.Range("L67").Value = 3624.9092
$ws.Range("M67").Value = -2634
$ws.Range("N67").Value = -5340.9092
$ws.Range("H93").Value = 87415
$ws.Range("J93").Value = 87415
$ws.Range("L93").Value = 87415
$ws.Range("N93").Value = -92407
$ws.Range("H98").Value = 1168.8636
$ws.Range("I98").Value = 1168.8636
$ws.Range("K98").Value = 1168.8636
$ws.Range("M98").Value = 329.1364000000001
$ws.Range("H106").Value = 3238.818
$ws.Range("I106").Value = 2223.5789
$ws.Range("K106").Value = 2223.5789
$ws.Range("M106").Value = -1592.5789
$ws.Range("H122").Value = 1168.8636
$ws.Range("I122").Value = 1168.8636
$ws.Range("K122").Value = 3506.5908
$ws.Range("M122").Value = -1056.5908

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 40016.25
$ws.Range("J24").Value = 40016.25
$ws.Range("L24").Value = 40016.25
$ws.Range("N24").Value = -40764.25
$ws.Range("H32").Value = 16953624
$ws.Range("I32").Value = 19234176
$ws.Range("J32").Value = 12382.571
$ws.Range("K32").Value = 19234176
$ws.Range("L32").Value = 12382.571
$ws.Range("M32").Value = -19233889
$ws.Range("N32").Value = -12956.571
$ws.Range("H63").Value = 3388.4614
$ws.Range("I63").Value = 2087.5
$ws.Range("J63").Value = 3966.6667
$ws.Range("K63").Value = 2087.5
$ws.Range("L63").Value = 3966.6667
$ws.Range("M63").Value = -1401.5
$ws.Range("N63").Value = -5338.6667
$ws.Range("H66").Value = 3388.4614
$ws.Range("I66").Value = 2087.5
$ws.Range("J66").Value = 3966.6667
$ws.Range("K66").Value = 10437.5
$ws.Range("L66").Value = 19833.3335
$ws.Range("M66").Value = -7005.5
$ws.Range("N66").Value = -26697.3335
$ws.Range("H74").Value = 2395
$ws.Range("I74").Value = 2651.375
$ws.Range("J74").Value = 1939.2222
$ws.Range("K74").Value = 2651.375
$ws.Range("L74").Value = 1939.2222
$ws.Range("M74").Value = -1777.375
$ws.Range("N74").Value = -3687.2222
$ws.Range("H77").Value = 2395
$ws.Range("I77").Value = 2651.375
$ws.Range("J77").Value = 1939.2222
$ws.Range("K77").Value = 13256.875
$ws.Range("L77").Value = 9696.110999999999
$ws.Range("M77").Value = -8888.875
$ws.Range("N77").Value = -18432.111
$ws.Range("H97").Value = 1180
$ws.Range("J97").Value = 1750
$ws.Range("L97").Value = 1750
$ws.Range("N97").Value = -2742
$ws.Range("H100").Value = 40016.25
$ws.Range("J100").Value = 40016.25
$ws.Range("L100").Value = 40016.25
$ws.Range("N100").Value = -42180.25
$ws.Range("H132").Value = 1549.7451
$ws.Range("I132").Value = 1410.6207
$ws.Range("K132").Value = 4231.8621
$ws.Range("M132").Value = -1701.8621

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 473.2
$ws.Range("I94").Value = 372.1111
$ws.Range("J94").Value = 624.8333
$ws.Range("K94").Value = 372.1111
$ws.Range("L94").Value = 624.8333
$ws.Range("M94").Value = 78.88889999999998
$ws.Range("N94").Value = -1526.8333
$ws.Range("H105").Value = 2601.4
$ws.Range("I105").Value = 1878.9656
$ws.Range("J105").Value = 2896.4788
$ws.Range("K105").Value = 1878.9656
$ws.Range("L105").Value = 2896.4788
$ws.Range("M105").Value = -131.9656
$ws.Range("N105").Value = -6390.4788
$ws.Range("H138").Value = 101640
$ws.Range("J138").Value = 101640
$ws.Range("L138").Value = 101640
$ws.Range("N138").Value = -111920
$ws.Range("H140").Value = 52000
$ws.Range("J140").Value = 52000
$ws.Range("L140").Value = 52000
$ws.Range("N140").Value = -62360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1828.6487
$ws.Range("I58").Value = 979.7241
$ws.Range("J58").Value = 4906
$ws.Range("K58").Value = 979.7241
$ws.Range("L58").Value = 4906
$ws.Range("M58").Value = -776.7241
$ws.Range("N58").Value = -5312
$ws.Range("H136").Value = 1828.6487
$ws.Range("I136").Value = 979.7241
$ws.Range("J136").Value = 4906
$ws.Range("K136").Value = 2939.1723
$ws.Range("L136").Value = 14718
$ws.Range("M136").Value = -389.1723000000002
$ws.Range("N136").Value = -19818

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 831.4
$ws.Range("I97").Value = 668.5
$ws.Range("J97").Value = 1483
$ws.Range("K97").Value = 668.5
$ws.Range("L97").Value = 1483
$ws.Range("M97").Value = -172.5
$ws.Range("N97").Value = -2475
$ws.Range("H132").Value = 3784.9387
$ws.Range("I132").Value = 4034.3713
$ws.Range("J132").Value = 3161.3572
$ws.Range("K132").Value = 12103.1139
$ws.Range("L132").Value = 9484.071599999999
$ws.Range("M132").Value = -9573.1139
$ws.Range("N132").Value = -14544.0716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 602.2759
$ws.Range("I55").Value = 535.16
$ws.Range("J55").Value = 1021.75
$ws.Range("K55").Value = 535.16
$ws.Range("L55").Value = 1021.75
$ws.Range("M55").Value = -362.16
$ws.Range("N55").Value = -1367.75
$ws.Range("H68").Value = 2526.5
$ws.Range("J68").Value = 2352.6
$ws.Range("L68").Value = 2352.6
$ws.Range("N68").Value = -3850.6
$ws.Range("H71").Value = 2526.5
$ws.Range("J71").Value = 2352.6
$ws.Range("L71").Value = 11763
$ws.Range("N71").Value = -19251
$ws.Range("H122").Value = 4139.8
$ws.Range("I122").Value = 3599
$ws.Range("J122").Value = 5401.6665
$ws.Range("K122").Value = 10797
$ws.Range("L122").Value = 16204.9995
$ws.Range("M122").Value = -8347
$ws.Range("N122").Value = -21104.9995
$ws.Range("H136").Value = 2104.8235
$ws.Range("I136").Value = 1568.2703
$ws.Range("J136").Value = 3522.8572
$ws.Range("K136").Value = 4704.810899999999
$ws.Range("L136").Value = 10568.5716
$ws.Range("M136").Value = -2154.810899999999
$ws.Range("N136").Value = -15668.5716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1583.7593
$ws.Range("I132").Value = 858.9706
$ws.Range("K132").Value = 2576.9118
$ws.Range("M132").Value = -46.91179999999986
